# Generate Report for Handback
# Update the timestamp values recorded on each localization status sheet.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for first data row.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-17 03:06:30"

# zh-cn sheet: handoff datetime (H2) and handback datetime (K2) for first data row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-17 03:06:25"
$wsZhCn.Range("K2").Value = "2016-08-17 03:06:42"

# de-de sheet: handback datetime (K2) for first data row.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-08-17 03:06:50"
